# Apply update to "International Ever Green_2024-12-24.xlsx":
#  - Append rows 32-41 to the "Orders" sheet (sheet1) with new flower line items.
#  - Extend the worksheet dimension / ignoredErrors range from L31 to L41 (handled
#    automatically by Excel once the new cells are populated).
#  - Update the "Summary" sheet (sheet2) G2 digest cell to include the new Number
#    values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Newline used inside a couple of the multi-line FlowerName values below.
$nl = "`n"

# row -> hashtable of column letter -> value (all stored as text, matching the
# existing sheet convention where every populated cell is a text string).
$rows = @(
    @{ Row = 32; C = "321_雪柳叶_Spiraea  leaves_undefined_1bunch"; F = "25" },
    @{ Row = 33; C = ("349_千层金绿_Melaleuca bracteata" + $nl + "（dyed orange）_Melaleuca bracteata F.Muell._1bunch"); F = "25" },
    @{ Row = 34; C = ("350_千层金红_Melaleuca bracteata" + $nl + "（dyed red）_Melaleuca bracteata F.Muell._1bunch"); F = "15" },
    @{ Row = 35; C = "109_绣球国产绿_Hydrangea Colombia Green (local)_Hydrangea L._1stem"; F = "40" },
    @{ Row = 36; C = "111_绣球单瓣紫粉_Hydrangea Purple&Pink S_Hydrangea L._1stem"; F = "60" },
    @{ Row = 37; C = "107_绣球单瓣浅粉_Hydrangea Light Pink S_Hydrangea L._1stem"; F = "50" },
    @{ Row = 38; A = "7"; C = "369_芦苇叶_undefined_undefined_1bunch"; F = "6" },
    @{ Row = 39; C = "688_山归来橙_undefined_undefined_1bunch"; F = "10" },
    @{ Row = 40; C = "595_玉兰叶_undefined_undefined_1bunch"; F = "13" },
    @{ Row = 41; C = "389_金合欢_mimosa_undefined_1bunch" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    if ($r.ContainsKey("A")) {
        $ws.Cells.Item($rowNum, 1).Value = "'" + $r.A
    }
    if ($r.ContainsKey("C")) {
        $ws.Cells.Item($rowNum, 3).Value = "'" + $r.C
    }
    if ($r.ContainsKey("F")) {
        $ws.Cells.Item($rowNum, 6).Value = "'" + $r.F
    }
}

# Update the Summary sheet's G2 running digest with the new Number (F) values
# appended, each as-is (missing/blank Number treated as "0").
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("G2").Value = "'" + "0101367310281538315101925184118101041019781253025252515406050610130"
